# Update gh-pages to output generated at 456a3b4
# Applies the data refresh to 展览 / 演出 / 全部类型 sheets:
#  - bump "want to go" (F) counts (and a couple of G price fixes) on existing rows
#  - insert one brand-new event row (COMIC WORLD) into 展览 (at row 26) and
#    全部类型 (at row 38), which pushes the remaining rows down by one
#  - the 火影only row picks up its own +1 bump after the shift

$wb = $excel.ActiveWorkbook

function Set-RowStyleA($ws, $addr) {
    # Reproduce cellXfs index 1 (bold font, thin box border, center/top align)
    # that the rest of column A in these tables uses.
    $ws.Range($addr).Font.Bold = $true
    $ws.Range($addr).HorizontalAlignment = -4108
    $ws.Range($addr).VerticalAlignment = -4160
    $ws.Range($addr).Borders.LineStyle = 1
}

# ---------------------------------------------------------------------
# Sheet "展览"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F4").Value2  = 1144
$ws1.Range("F7").Value2  = 605
$ws1.Range("F10").Value2 = 1446
$ws1.Range("F11").Value2 = 3066
$ws1.Range("F12").Value2 = 604
$ws1.Range("F13").Value2 = 1749
$ws1.Range("F14").Value2 = 1795
$ws1.Range("F15").Value2 = 843
$ws1.Range("F16").Value2 = 270
$ws1.Range("F17").Value2 = 1461
$ws1.Range("F18").Value2 = 286
$ws1.Range("F20").Value2 = 5
$ws1.Range("F21").Value2 = 1200
$ws1.Range("F22").Value2 = 397
$ws1.Range("F23").Value2 = 444
$ws1.Range("F24").Value2 = 91
$ws1.Range("F25").Value2 = 4735

# Insert the new event as row 26, pushing everything from the old row 26
# down to row 27 onward.
$ws1.Rows.Item(26).Insert()

$ws1.Range("A26").Value2 = 25
Set-RowStyleA $ws1 "A26"
$ws1.Range("B26").NumberFormat = "@"
$ws1.Range("B26").Value2 = "2024-06-08"
$ws1.Range("C26").Value2 = "广州·珠三角 2024 COMIC WORLD次元世界动漫游戏嘉年华"
$ws1.Range("D26").Value2 = "南洲路139号 小洲云文化艺术创意园"
$ws1.Range("E26").Value2 = "2024.06.08 10:00-06.10 17:00"
$ws1.Range("F26").Value2 = 1
$ws1.Range("G26").Value2 = 70
$ws1.Range("H26").Value2 = "https://show.bilibili.com/platform/detail.html?id=85020"
$ws1.Range("I26").Value2 = "//i2.hdslb.com/bfs/openplatform/202404/6g0jnqBP1714146665737.jpeg"

# The 火影only row (shifted from 30 to 31) gets its own +1 bump.
$ws1.Range("F31").Value2 = 118

# ---------------------------------------------------------------------
# Sheet "演出"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

$ws2.Range("F2").Value2  = 52
$ws2.Range("F3").Value2  = 27
$ws2.Range("F5").Value2  = 27
$ws2.Range("F7").Value2  = 63
$ws2.Range("G7").Value2  = 288
$ws2.Range("F13").Value2 = 25

# ---------------------------------------------------------------------
# Sheet "全部类型" (展览 + 演出 + 本地生活 combined)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F4").Value2  = 52
$ws4.Range("F5").Value2  = 27
$ws4.Range("F8").Value2  = 27
$ws4.Range("F10").Value2 = 63
$ws4.Range("G10").Value2 = 288
$ws4.Range("F12").Value2 = 1145
$ws4.Range("F15").Value2 = 605
$ws4.Range("F20").Value2 = 1446
$ws4.Range("F21").Value2 = 3066
$ws4.Range("F22").Value2 = 604
$ws4.Range("F23").Value2 = 1749
$ws4.Range("F24").Value2 = 1795
$ws4.Range("F25").Value2 = 843
$ws4.Range("F26").Value2 = 270
$ws4.Range("F27").Value2 = 1461
$ws4.Range("F28").Value2 = 286
$ws4.Range("F31").Value2 = 5
$ws4.Range("F33").Value2 = 1200
$ws4.Range("F34").Value2 = 397
$ws4.Range("F35").Value2 = 444
$ws4.Range("F36").Value2 = 91
$ws4.Range("F37").Value2 = 4735

# Insert the same new event as row 38.
$ws4.Rows.Item(38).Insert()

$ws4.Range("A38").Value2 = 37
Set-RowStyleA $ws4 "A38"
$ws4.Range("B38").NumberFormat = "@"
$ws4.Range("B38").Value2 = "2024-06-08"
$ws4.Range("C38").Value2 = "广州·珠三角 2024 COMIC WORLD次元世界动漫游戏嘉年华"
$ws4.Range("D38").Value2 = "南洲路139号 小洲云文化艺术创意园"
$ws4.Range("E38").Value2 = "2024.06.08 10:00-06.10 17:00"
$ws4.Range("F38").Value2 = 1
$ws4.Range("G38").Value2 = 70
$ws4.Range("H38").Value2 = "https://show.bilibili.com/platform/detail.html?id=85020"
$ws4.Range("I38").Value2 = "//i2.hdslb.com/bfs/openplatform/202404/6g0jnqBP1714146665737.jpeg"

# The 火影only row (shifted from 44 to 45) gets its own +1 bump.
$ws4.Range("F45").Value2 = 118
